$d = $word.ActiveDocument

# --- Prepend a run of 16 spaces before each of the three ingredient lines ---
# (Inserting at the very start of a multi-run paragraph creates a clean,
# separate leading run instead of merging into the existing text run.)
$d.Paragraphs.Item(2).Range.InsertBefore("                ")
$d.Paragraphs.Item(3).Range.InsertBefore("                ")
$d.Paragraphs.Item(4).Range.InsertBefore("                ")

# --- Rename the "DIRECTIONS" heading to "Instructions" ---
$found = $d.Content.Find.Execute("DIRECTIONS", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Instructions", 2)

# --- Add a new blank "\n" paragraph right after the grated-horseradish step ---
$gratedPara = $d.Paragraphs.Item(7)
$gratedPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item(8)
$newPara.Range.InsertBefore("\n")
